$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18 entry: date, hours worked, and description of the task.
$ws.Range("A18").Value = (Get-Date -Year 2021 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B18").Value = 3
$ws.Range("D18").Value = "Implementierung der Inputvalidierung der AcquirerOptions"

# Move the active selection to D19, mirroring where the cursor lands after
# finishing entry of the new row.
[void]$ws.Range("D19").Select()
